# Sync leve-profit sheets with latest scheduled market-board pull.
# Cell-level values below mirror the runner output for each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @(17, 8, 1322.4872),
    @(17, 10, 1368.5555),
    @(17, 12, 4105.666499999999),
    @(17, 14, -4441.666499999999),
    @(33, 8, 83),
    @(33, 9, 90.44444),
    @(33, 11, 90.44444),
    @(33, 13, 138.55556),
    @(70, 8, 180558820),
    @(70, 9, 125001250),
    @(70, 10, 208337580),
    @(70, 11, 375003750),
    @(70, 12, 625012740),
    @(70, 13, -375003480),
    @(70, 14, -625013280),
    @(73, 8, 180558820),
    @(73, 9, 125001250),
    @(73, 10, 208337580),
    @(73, 11, 375003750),
    @(73, 12, 625012740),
    @(73, 13, -375002814),
    @(73, 14, -625014612),
    @(74, 8, 33345312),
    @(74, 9, 71434310),
    @(74, 10, 17438),
    @(74, 11, 71434310),
    @(74, 12, 17438),
    @(74, 13, -71433374),
    @(74, 14, -19310),
    @(75, 8, 38821),
    @(75, 10, 45666.332),
    @(75, 12, 45666.332),
    @(75, 14, -47538.332),
    @(77, 8, 33345312),
    @(77, 9, 71434310),
    @(77, 10, 17438),
    @(77, 11, 357171550),
    @(77, 12, 87190),
    @(77, 13, -357166870),
    @(77, 14, -96550),
    @(78, 8, 38821),
    @(78, 10, 45666.332),
    @(78, 12, 136998.996),
    @(78, 14, -146358.996),
    @(87, 8, 49249.5),
    @(87, 10, 49249.5),
    @(87, 12, 49249.5),
    @(87, 14, -51745.5),
    @(90, 8, 49249.5),
    @(90, 10, 49249.5),
    @(90, 12, 147748.5),
    @(90, 14, -160228.5),
    @(98, 8, 2242.879),
    @(98, 9, 2258.4194),
    @(98, 10, 2002),
    @(98, 11, 2258.4194),
    @(98, 12, 2002),
    @(98, 13, -760.4194000000002),
    @(98, 14, -4998),
    @(113, 8, 27791720),
    @(113, 9, 2538),
    @(113, 10, 41686308),
    @(113, 11, 2538),
    @(113, 12, 41686308),
    @(113, 13, 716),
    @(113, 14, -41692816),
    @(115, 8, 959.6667),
    @(115, 9, 992),
    @(115, 10, 895),
    @(115, 11, 2976),
    @(115, 12, 2685),
    @(115, 13, -1409),
    @(115, 14, -5819),
    @(122, 8, 2242.879),
    @(122, 9, 2258.4194),
    @(122, 10, 2002),
    @(122, 11, 6775.2582),
    @(122, 12, 6006),
    @(122, 13, -4325.2582),
    @(122, 14, -10906),
    @(131, 8, 1887.8823),
    @(131, 9, 1506.9286),
    @(131, 11, 4520.7858),
    @(131, 13, 519.2142000000003),
    @(135, 8, 667351.2),
    @(135, 9, 769748.6),
    @(135, 11, 6927737.399999999),
    @(135, 13, -6925202.399999999),
    @(138, 8, 1337401.5),
    @(138, 9, 2596),
    @(138, 10, 1856492.5),
    @(138, 11, 7788),
    @(138, 12, 5569477.5),
    @(138, 13, -2648),
    @(138, 14, -5579757.5),
    @(139, 8, 0),
    @(139, 10, 0),
    @(139, 12, 0),
    @(139, 14, $null),
    @(141, 8, 6552.909),
    @(141, 9, 6012.5713),
    @(141, 11, 18037.7139),
    @(141, 13, -12857.7139)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @(2, 8, 1534.0333),
    @(2, 9, 1307.1),
    @(2, 10, 1987.9),
    @(2, 11, 1307.1),
    @(2, 12, 1987.9),
    @(2, 13, -1194.1),
    @(2, 14, -2213.9),
    @(32, 8, 3181741.8),
    @(32, 9, 3395249.8),
    @(32, 11, 3395249.8),
    @(32, 13, -3394962.8),
    @(45, 8, 4569.8647),
    @(45, 9, 3231.8667),
    @(45, 11, 3231.8667),
    @(45, 13, -2854.8667),
    @(57, 8, 5499.364),
    @(57, 9, 5499.364),
    @(57, 11, 5499.364),
    @(57, 13, -5015.364),
    @(88, 8, 1760.1),
    @(88, 10, 1760.1),
    @(88, 12, 1760.1),
    @(88, 14, -2572.1),
    @(91, 8, 1760.1),
    @(91, 10, 1760.1),
    @(91, 12, 1760.1),
    @(91, 14, -4568.1),
    @(97, 8, 3624695),
    @(97, 9, 1246),
    @(97, 11, 1246),
    @(97, 13, -750),
    @(116, 8, 1534.0333),
    @(116, 9, 1307.1),
    @(116, 10, 1987.9),
    @(116, 11, 1307.1),
    @(116, 12, 1987.9),
    @(116, 13, 986.9000000000001),
    @(116, 14, -6575.9),
    @(122, 8, 3621.7273),
    @(122, 9, 2945.2222),
    @(122, 11, 8835.6666),
    @(122, 13, -6385.6666),
    @(126, 8, 5397.4443),
    @(126, 9, 5397.4443),
    @(126, 11, 16192.3329),
    @(126, 13, -13722.3329),
    @(129, 8, 86000),
    @(129, 10, 86000),
    @(129, 12, 86000),
    @(129, 14, -96000),
    @(132, 8, 3256.7754),
    @(132, 9, 1068),
    @(132, 10, 7771.125),
    @(132, 11, 3204),
    @(132, 12, 23313.375),
    @(132, 13, -674),
    @(132, 14, -28373.375),
    @(140, 8, 0),
    @(140, 10, 0),
    @(140, 12, 0),
    @(140, 14, $null),
    @(141, 8, 44713.5),
    @(141, 10, 44713.5),
    @(141, 12, 44713.5),
    @(141, 14, -55073.5)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @(3, 8, 1534.0333),
    @(3, 9, 1307.1),
    @(3, 10, 1987.9),
    @(3, 11, 1307.1),
    @(3, 12, 1987.9),
    @(3, 13, -1193.1),
    @(3, 14, -2215.9),
    @(86, 8, 13947145),
    @(86, 9, 19309970),
    @(86, 10, 3799.2),
    @(86, 11, 19309970),
    @(86, 12, 3799.2),
    @(86, 13, -19308847),
    @(86, 14, -6045.2),
    @(89, 8, 13947145),
    @(89, 9, 19309970),
    @(89, 10, 3799.2),
    @(89, 11, 96549850),
    @(89, 12, 18996),
    @(89, 13, -96544234),
    @(89, 14, -30228),
    @(94, 8, 2450.1052),
    @(94, 9, 1447.7142),
    @(94, 10, 5256.8),
    @(94, 11, 1447.7142),
    @(94, 12, 5256.8),
    @(94, 13, -996.7141999999999),
    @(94, 14, -6158.8),
    @(105, 8, 3379.2144),
    @(105, 9, 2780.65),
    @(105, 11, 2780.65),
    @(105, 13, -1033.65),
    @(134, 8, 7817374),
    @(134, 9, 11908108),
    @(134, 11, 35724324),
    @(134, 13, -35721789)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @(16, 8, 6452.8887),
    @(16, 9, 1315.8),
    @(16, 10, 8428.691999999999),
    @(16, 11, 1315.8),
    @(16, 12, 8428.691999999999),
    @(16, 13, -1028.8),
    @(16, 14, -9002.691999999999),
    @(69, 8, 27616.727),
    @(69, 9, 24723),
    @(69, 11, 24723),
    @(69, 13, -23974),
    @(72, 8, 27616.727),
    @(72, 9, 24723),
    @(72, 11, 74169),
    @(72, 13, -70425),
    @(76, 8, 5299.125),
    @(76, 9, 5299.125),
    @(76, 11, 5299.125),
    @(76, 13, -4984.125),
    @(79, 8, 5299.125),
    @(79, 9, 5299.125),
    @(79, 11, 5299.125),
    @(79, 13, -4207.125),
    @(86, 8, 18387628),
    @(86, 9, 26047570),
    @(86, 10, 3764.8),
    @(86, 11, 26047570),
    @(86, 12, 3764.8),
    @(86, 13, -26046447),
    @(86, 14, -6010.8),
    @(89, 8, 18387628),
    @(89, 9, 26047570),
    @(89, 10, 3764.8),
    @(89, 11, 130237850),
    @(89, 12, 18824),
    @(89, 13, -130232234),
    @(89, 14, -30056),
    @(94, 8, 1711.6666),
    @(94, 9, 2136.4),
    @(94, 11, 2136.4),
    @(94, 13, -1685.4),
    @(99, 8, 3124),
    @(99, 9, 2148.8),
    @(99, 11, 2148.8),
    @(99, 13, -650.8000000000002),
    @(113, 8, 6452.8887),
    @(113, 9, 1315.8),
    @(113, 10, 8428.691999999999),
    @(113, 11, 1315.8),
    @(113, 12, 8428.691999999999),
    @(113, 13, 854.2),
    @(113, 14, -12768.692),
    @(126, 8, 3124),
    @(126, 9, 2148.8),
    @(126, 11, 6446.400000000001),
    @(126, 13, -3976.400000000001),
    @(133, 8, 0),
    @(133, 10, 0),
    @(133, 12, 0),
    @(133, 14, $null),
    @(134, 8, 2707.5908),
    @(134, 9, 1475.6364),
    @(134, 10, 6403.4546),
    @(134, 11, 4426.9092),
    @(134, 12, 19210.3638),
    @(134, 13, -1891.9092),
    @(134, 14, -24280.3638)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @(68, 8, 33338846),
    @(68, 10, 10107.667),
    @(68, 12, 30323.001),
    @(68, 14, -31945.001),
    @(71, 8, 33338846),
    @(71, 10, 10107.667),
    @(71, 12, 90969.003),
    @(71, 14, -99081.003),
    @(92, 8, 19232468),
    @(92, 9, 1899.5),
    @(92, 11, 5698.5),
    @(92, 13, -4450.5),
    @(97, 8, 358.5),
    @(97, 10, 387),
    @(97, 12, 1161),
    @(97, 14, -2153),
    @(98, 8, 2247.6667),
    @(98, 10, 2296.6),
    @(98, 12, 6889.799999999999),
    @(98, 14, -9885.799999999999),
    @(129, 8, 223522.89),
    @(129, 9, 881.8),
    @(129, 10, 501824.25),
    @(129, 11, 2645.4),
    @(129, 12, 1505472.75),
    @(129, 13, 2354.6),
    @(129, 14, -1515472.75)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @(70, 8, 18682.945),
    @(70, 9, 24087.9),
    @(70, 10, 11926.75),
    @(70, 11, 24087.9),
    @(70, 12, 11926.75),
    @(70, 13, -23817.9),
    @(70, 14, -12466.75),
    @(73, 8, 18682.945),
    @(73, 9, 24087.9),
    @(73, 10, 11926.75),
    @(73, 11, 24087.9),
    @(73, 12, 11926.75),
    @(73, 13, -23151.9),
    @(73, 14, -13798.75),
    @(102, 8, 6362.3),
    @(102, 9, 6199.5),
    @(102, 11, 6199.5),
    @(102, 13, -4577.5),
    @(122, 8, 4542294.5),
    @(122, 9, 5588989),
    @(122, 10, 6619.3335),
    @(122, 11, 16766967),
    @(122, 12, 19858.0005),
    @(122, 13, -16764517),
    @(122, 14, -24758.0005)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @(7, 8, 5714.3335),
    @(7, 9, 5191.75),
    @(7, 10, 5975.625),
    @(7, 11, 5191.75),
    @(7, 12, 5975.625),
    @(7, 13, -5079.75),
    @(7, 14, -6199.625),
    @(22, 8, 2250.3),
    @(22, 9, 856),
    @(22, 10, 5503.6665),
    @(22, 11, 856),
    @(22, 12, 5503.6665),
    @(22, 13, -561),
    @(22, 14, -6093.6665),
    @(27, 8, 2250.3),
    @(27, 9, 856),
    @(27, 10, 5503.6665),
    @(27, 11, 856),
    @(27, 12, 5503.6665),
    @(27, 13, -749),
    @(27, 14, -5717.6665),
    @(40, 8, 5728.294),
    @(40, 9, 4609.778),
    @(40, 10, 6986.625),
    @(40, 11, 4609.778),
    @(40, 12, 6986.625),
    @(40, 13, -4473.778),
    @(40, 14, -7258.625),
    @(46, 8, 1719.3158),
    @(46, 9, 1170.3125),
    @(46, 10, 4647.3335),
    @(46, 11, 1170.3125),
    @(46, 12, 4647.3335),
    @(46, 13, -982.3125),
    @(46, 14, -5023.3335),
    @(61, 8, 7705.1333),
    @(61, 9, 7051.1665),
    @(61, 11, 7051.1665),
    @(61, 13, -6849.1665),
    @(107, 8, 4561.625),
    @(107, 9, 4561.625),
    @(107, 11, 4561.625),
    @(107, 13, -2641.625),
    @(113, 8, 7705.1333),
    @(113, 9, 7051.1665),
    @(113, 11, 7051.1665),
    @(113, 13, -4881.1665),
    @(122, 8, 3497.568),
    @(122, 9, 2801.4),
    @(122, 10, 4989.357),
    @(122, 11, 8404.200000000001),
    @(122, 12, 14968.071),
    @(122, 13, -5954.200000000001),
    @(122, 14, -19868.071),
    @(126, 8, 5714.3335),
    @(126, 9, 5191.75),
    @(126, 10, 5975.625),
    @(126, 11, 15575.25),
    @(126, 12, 17926.875),
    @(126, 13, -13105.25),
    @(126, 14, -22866.875),
    @(132, 8, 9619433),
    @(132, 9, 16131243),
    @(132, 11, 48393729),
    @(132, 13, -48391199),
    @(136, 8, 8856.633),
    @(136, 9, 2481.9656),
    @(136, 10, 18099.9),
    @(136, 11, 7445.8968),
    @(136, 12, 54299.7),
    @(136, 13, -4895.8968),
    @(136, 14, -59399.7),
    @(140, 8, 79833.336),
    @(140, 10, 79833.336),
    @(140, 12, 79833.336),
    @(140, 14, -90193.336)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @(107, 8, 1181.6957),
    @(107, 9, 477),
    @(107, 10, 1950.4546),
    @(107, 11, 1431),
    @(107, 12, 5851.3638),
    @(107, 13, 489),
    @(107, 14, -9691.363799999999),
    @(113, 8, 1438.68),
    @(113, 9, 1252.8462),
    @(113, 10, 1640),
    @(113, 11, 3758.5386),
    @(113, 12, 4920),
    @(113, 13, -1588.5386),
    @(113, 14, -9260),
    @(122, 8, 576341.3),
    @(122, 9, 1002401),
    @(122, 11, 3007203),
    @(122, 13, -3004753),
    @(126, 8, 746),
    @(126, 9, 746),
    @(126, 11, 2238),
    @(126, 13, 232),
    @(133, 8, 176250),
    @(133, 10, 176250),
    @(133, 12, 176250),
    @(133, 14, -186370)
)
foreach ($u in $updates) {
    if ($u[2] -eq $null) {
        $ws.Cells.Item($u[0], $u[1]).ClearContents()
    } else {
        $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
    }
}
